$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value  = "2.0.0"
$meta.Range("B8").Value  = "2024-06-14T10:04:53+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"
$meta.Range("B11").Value = "Intervention codes allowed in this implementation guide."

# --- Rename the second sheet -------------------------------------------------
$inc = $wb.Worksheets.Item("Include ValueSets")
$inc.Name = "Include from FSIII"

# --- Replace the "Include from FSIII" sheet content -------------------------
# Header row
$inc.Cells.Item(1, 1).Value = "Concept"
$inc.Cells.Item(1, 2).Value = "Description"

$concepts = @(
    "01a500f6-c221-4fd0-b518-cd361218b60d",
    "03a3ebdb-9e2d-4be1-b32b-42f0bd2a3362",
    "61692d91-69b8-4830-9453-3d58454e49d3",
    "6d24992e-e0a2-43e7-bc27-0234622a8655",
    "6eddbaf7-2a73-49d4-91e7-6138d419f58c",
    "924e9828-84cf-4689-9551-0ebb6dc71b98",
    "ab87c0b5-40be-4e0a-b749-d9f833bfed2d",
    "abe847e0-1ce0-44dc-a675-ce05b66f47e6",
    "c9a99304-1788-43b7-b7be-e187b092ae9c",
    "cf7a55c2-7061-47ed-b7c5-e29620fe93bf",
    "d1e016b5-150a-4ac4-97ba-d3e19e28471e",
    "e71b7d85-5c78-49c2-8624-8499d162317b",
    "ee5606ac-1bed-487e-aa3c-72dcc30ec037",
    "f30cab6d-2a42-4358-99d7-811127fb6e05"
)

$row = 2
foreach ($concept in $concepts) {
    $inc.Cells.Item($row, 1).Value = $concept
    $row = $row + 1
}

# Row 16: blank concept / blank description
$inc.Cells.Item(16, 1).Value = ""
$inc.Cells.Item(16, 2).Value = ""

# Row 17: System URI
$inc.Cells.Item(17, 1).Value = "System URI"
$inc.Cells.Item(17, 2).Value = "urn:oid:1.2.208.176.2.21"
